# Apply market-price / profit recalculations to the Leve sheets
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 2893.3
$ws.Range("J32").Value = 3022.5
$ws.Range("L32").Value = 3022.5
$ws.Range("N32").Value = -3674.5
# Row 33
$ws.Range("H33").Value = 273.36365
$ws.Range("I33").Value = 185.25
$ws.Range("J33").Value = 508.33334
$ws.Range("K33").Value = 185.25
$ws.Range("L33").Value = 508.33334
$ws.Range("M33").Value = 43.75
$ws.Range("N33").Value = -966.33334
# Row 42
$ws.Range("H42").Value = 13908
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 13908
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 41724
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -42184
# Row 47
$ws.Range("H47").Value = 13754
$ws.Range("I47").Value = 10005.333
$ws.Range("K47").Value = 10005.333
$ws.Range("M47").Value = -9033.333000000001
# Row 132
$ws.Range("H132").Value = 1923.1666
$ws.Range("I132").Value = 1923.1666
$ws.Range("K132").Value = 5769.4998
$ws.Range("M132").Value = -3239.4998
# Row 137
$ws.Range("H137").Value = 744.7857
$ws.Range("I137").Value = 428
$ws.Range("J137").Value = 797.5833
$ws.Range("K137").Value = 1284
$ws.Range("L137").Value = 2392.7499
$ws.Range("M137").Value = 1266
$ws.Range("N137").Value = -7492.7499
# Row 138
$ws.Range("H138").Value = 3979.9
$ws.Range("J138").Value = 4405.4136
$ws.Range("L138").Value = 13216.2408
$ws.Range("N138").Value = -23496.2408

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1403.5
$ws.Range("I2").Value = 855.5
$ws.Range("K2").Value = 855.5
$ws.Range("M2").Value = -742.5
# Row 4
$ws.Range("H4").Value = 2331.8333
$ws.Range("I4").Value = 1997.6666
$ws.Range("J4").Value = 2666
$ws.Range("K4").Value = 1997.6666
$ws.Range("L4").Value = 2666
$ws.Range("M4").Value = -1881.6666
$ws.Range("N4").Value = -2898
# Row 74
$ws.Range("H74").Value = 1048.1666
$ws.Range("I74").Value = 947.25
$ws.Range("K74").Value = 947.25
$ws.Range("M74").Value = -73.25
# Row 77
$ws.Range("H77").Value = 1048.1666
$ws.Range("I77").Value = 947.25
$ws.Range("K77").Value = 4736.25
$ws.Range("M77").Value = -368.25
# Row 116
$ws.Range("H116").Value = 1403.5
$ws.Range("I116").Value = 855.5
$ws.Range("K116").Value = 855.5
$ws.Range("M116").Value = 1438.5

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1403.5
$ws.Range("I3").Value = 855.5
$ws.Range("K3").Value = 855.5
$ws.Range("M3").Value = -741.5
# Row 86
$ws.Range("H86").Value = 7901
$ws.Range("I86").Value = 8035
$ws.Range("K86").Value = 8035
$ws.Range("M86").Value = -6912
# Row 89
$ws.Range("H89").Value = 7901
$ws.Range("I89").Value = 8035
$ws.Range("K89").Value = 40175
$ws.Range("M89").Value = -34559
# Row 105
$ws.Range("H105").Value = 2052.2
$ws.Range("I105").Value = 1712.5
$ws.Range("K105").Value = 1712.5
$ws.Range("M105").Value = 34.5
# Row 134
$ws.Range("H134").Value = 3480.1765
$ws.Range("I134").Value = 3572.6875
$ws.Range("K134").Value = 10718.0625
$ws.Range("M134").Value = -8183.0625

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 103.8421
$ws.Range("I7").Value = 64.13333
$ws.Range("J7").Value = 252.75
$ws.Range("K7").Value = 64.13333
$ws.Range("L7").Value = 252.75
$ws.Range("M7").Value = 48.86667
$ws.Range("N7").Value = -478.75
# Row 31
$ws.Range("H31").Value = 1236.7693
$ws.Range("I31").Value = 1152.6666
$ws.Range("J31").Value = 1281.2941
$ws.Range("K31").Value = 1152.6666
$ws.Range("L31").Value = 1281.2941
$ws.Range("M31").Value = -857.6666
$ws.Range("N31").Value = -1871.2941
# Row 34
$ws.Range("H34").Value = 1236.7693
$ws.Range("I34").Value = 1152.6666
$ws.Range("J34").Value = 1281.2941
$ws.Range("K34").Value = 1152.6666
$ws.Range("L34").Value = 1281.2941
$ws.Range("M34").Value = -950.6666
$ws.Range("N34").Value = -1685.2941

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 1337
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 4500
$ws.Range("N22").Value = -4838
# Row 27
$ws.Range("H27").Value = 1337
$ws.Range("J27").Value = 1500
$ws.Range("L27").Value = 4500
$ws.Range("N27").Value = -4704
# Row 88
$ws.Range("H88").Value = 17665.334
$ws.Range("J88").Value = 17665.334
$ws.Range("L88").Value = 52996.00199999999
$ws.Range("N88").Value = -53852.00199999999
# Row 91
$ws.Range("H91").Value = 17665.334
$ws.Range("J91").Value = 17665.334
$ws.Range("L91").Value = 52996.00199999999
$ws.Range("N91").Value = -55960.00199999999
# Row 95
$ws.Range("H95").Value = 8889.333000000001
$ws.Range("J95").Value = 8889.333000000001
$ws.Range("L95").Value = 26667.999
$ws.Range("N95").Value = -30785.999

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2563.889
$ws.Range("I102").Value = 1994.5
$ws.Range("K102").Value = 1994.5
$ws.Range("M102").Value = -372.5
# Row 126
$ws.Range("H126").Value = 5592.7
$ws.Range("I126").Value = 4450
$ws.Range("K126").Value = 13350
$ws.Range("M126").Value = -10880

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3499.25
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 3999
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 3999
$ws.Range("M22").Value = -1705
$ws.Range("N22").Value = -4589
# Row 27
$ws.Range("H27").Value = 3499.25
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 3999
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 3999
$ws.Range("M27").Value = -1893
$ws.Range("N27").Value = -4213
# Row 46
$ws.Range("H46").Value = 1914.4
$ws.Range("I46").Value = 2469
$ws.Range("J46").Value = 1359.8
$ws.Range("K46").Value = 2469
$ws.Range("L46").Value = 1359.8
$ws.Range("M46").Value = -2281
$ws.Range("N46").Value = -1735.8
# Row 55
$ws.Range("H55").Value = 1972.4
$ws.Range("I55").Value = 1689.6666
$ws.Range("J55").Value = 2396.5
$ws.Range("K55").Value = 1689.6666
$ws.Range("L55").Value = 2396.5
$ws.Range("M55").Value = -1516.6666
$ws.Range("N55").Value = -2742.5
# Row 100
$ws.Range("H100").Value = 1421.8572
$ws.Range("I100").Value = 1408.8334
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 1408.8334
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -867.8334
$ws.Range("N100").Value = -2582

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
# Row 11
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
# Row 96
$ws.Range("H96").Value = 4795
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
